$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.243.36"
$ws.Range("D3").Value = "1.906.17"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "326.30"
$ws.Range("E5").Value = "  -0.43%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  -0.20%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4651"
$ws.Range("E7").Value = "  +0.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3921"
$ws.Range("E8").Value = "  -0.84%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07898"
$ws.Range("E9").Value = "  -0.82%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9923"
$ws.Range("E10").Value = "  -0.90%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "21.86"
$ws.Range("E11").Value = "  -1.98%  "
$ws.Range("D12").Value = "1.915.00"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "7.084"
$ws.Range("E13").Value = "  -0.62%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.749"
$ws.Range("E14").Value = "  -0.20%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.06991"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "88.27"
$ws.Range("E16").Value = "  -0.57%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.002"
$ws.Range("E17").Value = "  -0.22%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001000"
$ws.Range("E18").Value = "  -0.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.14"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("E20").Value = "  -0.13%  "
$ws.Range("D21").Value = "29.244.60"
$ws.Range("E21").Value = "  -0.13%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.317"
$ws.Range("E22").Value = "  -0.80%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.10"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("B24").Value = "Toncoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.109"
$ws.Range("E24").Value = "  +2.86%  "
$ws.Range("B25").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C25").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D25").Value = "2.095.82"
$ws.Range("E25").Value = "  -4.14%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "156.25"
$ws.Range("E26").Value = "  -0.30%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.43"
$ws.Range("E27").Value = "  -0.38%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.978"
$ws.Range("E28").Value = "  +0.86%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "118.74"
$ws.Range("E29").Value = "  -0.26%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.890"
$ws.Range("E30").Value = "  -5.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09357"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9036"
$ws.Range("E32").Value = "  -2.46%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.261"
$ws.Range("E33").Value = "  -1.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.330"
$ws.Range("E34").Value = "  -1.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.181"
$ws.Range("E35").Value = "  -2.44%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.182"
$ws.Range("E36").Value = "  +1.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.05782"
$ws.Range("E37").Value = "  -0.94%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02093"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.714"
$ws.Range("E40").Value = "  -3.61%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5713"
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1794"
$ws.Range("E42").Value = "  -0.62%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "9.743"
$ws.Range("E43").Value = "  -2.44%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "11.97"
$ws.Range("E44").Value = "  -0.51%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.5366"
$ws.Range("E45").Value = "  -1.18%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.182"
$ws.Range("E46").Value = "  -1.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.07019"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.859"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.571"
$ws.Range("E49").Value = "  -0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "113.46"
$ws.Range("E50").Value = "  +1.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.050"
$ws.Range("E51").Value = "  -1.18%  "
